# Applies updated crypto price/volume figures to Sheet1 (D = Price, E = Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.848.03'
$ws.Range("E2").Value = '  -5.75%  '
$ws.Range("D3").Value = '2.582.70'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.70'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.70'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.581'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.44%  '
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.563'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.00'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.72%  '
$ws.Range("E11").Value = '  -3.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.83'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.54%  '
$ws.Range("D13").Value = '2.977.86'
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("E14").Value = '  +1.18%  '
$ws.Range("D15").Value = '2.584.63'
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.895'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.39'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.05%  '
$ws.Range("D18").Value = '43.804.99'
$ws.Range("E18").Value = '  -6.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.72'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.67'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '266.30'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.23'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.94'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.47'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.25%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.26'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.59%  '
$ws.Range("E29").Value = '  -3.57%  '
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.58'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.21'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.66'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.80'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0818'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.93%  '
$ws.Range("E37").Value = '  -4.40%  '
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.33'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.95'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.59'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0317'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.88'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.13%  '
$ws.Range("D44").Value = '2.039.01'
$ws.Range("E44").Value = '  -4.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.16'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.17'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.99%  '
$ws.Range("E48").Value = '  +4.56%  '
$ws.Range("D49").Value = '2.835.93'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '105.71'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.192'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.36%  '
